# Updates NATMI ligand-receptor pair TPM statistics for rows 2-13 (columns E:T)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @{ Row = 2; Values = @{ "E" = "3"; "F" = "1"; "G" = "8.223245666666665"; "H" = "24.669737"; "I" = "0.003010099901484358"; "J" = "0.003010099901484359"; "K" = "2"; "L" = "0.6666666666666666"; "M" = "1.417175333333333"; "N" = "4.251526"; "O" = "0.02388929759887871"; "P" = "0.02388929759887871"; "Q" = "11.65378091874022"; "R" = "104.884028268662"; "S" = "7.19091723489153E-05"; "T" = "7.190917234891532E-05" } },
    @{ Row = 3; Values = @{ "E" = "3"; "F" = "1"; "G" = "8.223245666666665"; "H" = "24.669737"; "I" = "0.003010099901484358"; "J" = "0.003010099901484359"; "K" = "3"; "L" = "1"; "M" = "2.226320666666667"; "N" = "6.678962"; "O" = "0.0375290450698413"; "P" = "0.03752904506984131"; "Q" = "18.30758177477711"; "R" = "164.768235972994"; "S" = "0.0001129661748675313"; "T" = "0.0001129661748675314" } },
    @{ Row = 4; Values = @{ "E" = "3"; "F" = "1"; "G" = "8.223245666666665"; "H" = "24.669737"; "I" = "0.003010099901484358"; "J" = "0.003010099901484359"; "K" = "3"; "L" = "1"; "M" = "55.67910766666667"; "N" = "167.037323"; "O" = "0.93858165733128"; "P" = "0.93858165733128"; "Q" = "457.8629808437834"; "R" = "4120.766827594051"; "S" = "0.002825224554267911"; "T" = "0.002825224554267912" } },
    @{ Row = 5; Values = @{ "E" = "3"; "F" = "1"; "G" = "2706.934895666667"; "H" = "8120.804687"; "I" = "0.9908672065823976"; "J" = "0.9908672065823977"; "K" = "2"; "L" = "0.6666666666666666"; "M" = "1.417175333333333"; "N" = "4.251526"; "O" = "0.02388929759887871"; "P" = "0.02388929759887871"; "Q" = "3836.201363078041"; "R" = "34525.81226770236"; "S" = "0.02367112157901652"; "T" = "0.02367112157901653" } },
    @{ Row = 6; Values = @{ "E" = "3"; "F" = "1"; "G" = "2706.934895666667"; "H" = "8120.804687"; "I" = "0.9908672065823976"; "J" = "0.9908672065823977"; "K" = "3"; "L" = "1"; "M" = "2.226320666666667"; "N" = "6.678962"; "O" = "0.0375290450698413"; "P" = "0.03752904506984131"; "Q" = "6026.505101543878"; "R" = "54238.54591389489"; "S" = "0.03718630005405855"; "T" = "0.03718630005405856" } },
    @{ Row = 7; Values = @{ "E" = "3"; "F" = "1"; "G" = "2706.934895666667"; "H" = "8120.804687"; "I" = "0.9908672065823976"; "J" = "0.9908672065823977"; "K" = "3"; "L" = "1"; "M" = "55.67910766666667"; "N" = "167.037323"; "O" = "0.93858165733128"; "P" = "0.93858165733128"; "Q" = "150719.7195024815"; "R" = "1356477.475522333"; "S" = "0.9300097849493225"; "T" = "0.9300097849493226" } },
    @{ Row = 8; Values = @{ "E" = "3"; "F" = "1"; "G" = "14.14340733333333"; "H" = "42.430222"; "I" = "0.005177161275053701"; "J" = "0.005177161275053702"; "K" = "2"; "L" = "0.6666666666666666"; "M" = "1.417175333333333"; "N" = "4.251526"; "O" = "0.02388929759887871"; "P" = "0.02388929759887871"; "Q" = "20.04368800208578"; "R" = "180.393192018772"; "S" = "0.0001236787464171482"; "T" = "0.0001236787464171482" } },
    @{ Row = 9; Values = @{ "E" = "3"; "F" = "1"; "G" = "14.14340733333333"; "H" = "42.430222"; "I" = "0.005177161275053701"; "J" = "0.005177161275053702"; "K" = "3"; "L" = "1"; "M" = "2.226320666666667"; "N" = "6.678962"; "O" = "0.0375290450698413"; "P" = "0.03752904506984131"; "Q" = "31.48776004328489"; "R" = "283.389840389564"; "S" = "0.0001942939188253274"; "T" = "0.0001942939188253275" } },
    @{ Row = 10; Values = @{ "E" = "3"; "F" = "1"; "G" = "14.14340733333333"; "H" = "42.430222"; "I" = "0.005177161275053701"; "J" = "0.005177161275053702"; "K" = "3"; "L" = "1"; "M" = "55.67910766666667"; "N" = "167.037323"; "O" = "0.93858165733128"; "P" = "0.93858165733128"; "Q" = "787.4922996861897"; "R" = "7087.430697175707"; "S" = "0.004859188609811225"; "T" = "0.004859188609811226" } },
    @{ Row = 11; Values = @{ "E" = "3"; "F" = "1"; "G" = "2.583085"; "H" = "7.749255"; "I" = "0.0009455322410643118"; "J" = "0.0009455322410643119"; "K" = "2"; "L" = "0.6666666666666666"; "M" = "1.417175333333333"; "N" = "4.251526"; "O" = "0.02388929759887871"; "P" = "0.02388929759887871"; "Q" = "3.660684345903334"; "R" = "32.94615911313"; "S" = "2.258810109612007E-05"; "T" = "2.258810109612007E-05" } },
    @{ Row = 12; Values = @{ "E" = "3"; "F" = "1"; "G" = "2.583085"; "H" = "7.749255"; "I" = "0.0009455322410643118"; "J" = "0.0009455322410643119"; "K" = "3"; "L" = "1"; "M" = "2.226320666666667"; "N" = "6.678962"; "O" = "0.0375290450698413"; "P" = "0.03752904506984131"; "Q" = "5.750775519256667"; "R" = "51.75697967331"; "S" = "3.548492208989061E-05"; "T" = "3.548492208989062E-05" } },
    @{ Row = 13; Values = @{ "E" = "3"; "F" = "1"; "G" = "2.583085"; "H" = "7.749255"; "I" = "0.0009455322410643118"; "J" = "0.0009455322410643119"; "K" = "3"; "L" = "1"; "M" = "55.67910766666667"; "N" = "167.037323"; "O" = "0.93858165733128"; "P" = "0.93858165733128"; "Q" = "143.8238678271517"; "R" = "1294.414810444365"; "S" = "0.0008874592178783011"; "T" = "0.0008874592178783012" } }
)

foreach ($rowEntry in $rowData) {
    $rowNum = $rowEntry.Row
    foreach ($col in $rowEntry.Values.Keys) {
        $addr = "$col$rowNum"
        $ws.Range($addr).Value = [double]$rowEntry.Values[$col]
    }
}
